$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59: drop trailing "!" from rejected-subject translations
$ws.Range("B59").Value = "{0}, Tú solicitud de emisión de certificado {1} ha sido rechazada"
$ws.Range("C59").Value = "{0}, Your request to issue a certificate {1} has been rejected"

# New notification rows 85-103 (FCM gateway notification strings)
$ws.Range("A85").Value = "notification_certificate_enabled_title"
$ws.Range("B85").Value = "{0}, Tú certificado {1} en Trust Certification System ha sido habilitado."
$ws.Range("C85").Value = "{0}, Your certificate {1} in the Trust Certification System has been enabled."
$ws.Range("C85").Font.Size = 11

$ws.Range("A86").Value = "notification_certificate_enabled_body"
$ws.Range("B86").Value = "{0}, Tú certificado {1} ha sido habilitado, podrás renovar y utilizar tu certificado con normalidad."
$ws.Range("C86").Value = "{0}, Your certificate {1} has been enabled, you can renew and use your certificate normally."
$ws.Range("C86").Font.Size = 11

$ws.Range("A87").Value = "notification_certificate_disabled_title"
$ws.Range("B87").Value = "{0}, Tú certificado {1} en Trust Certification System ha sido deshabilitado."
$ws.Range("C87").Value = "{0}, Your certificate {1} in the Trust Certification System has been disabled."
$ws.Range("C87").Font.Size = 11

$ws.Range("A88").Value = "notification_certificate_disabled_body"
$ws.Range("B88").Value = "{0}, Tú certificado {1} ha sido deshabilitado, este no podrá ser renovado ni utilizado hasta que sea habilitado de nuevo."
$ws.Range("C88").Value = "{0}, Your certificate {1} has been disabled, it cannot be renewed or used until it is enabled again."
$ws.Range("C88").Font.Size = 11

$ws.Range("A89").Value = "notification_certificate_issued_title"
$ws.Range("B89").Value = "{0}, tú certificado {1} ha sido registrado con éxito!"
$ws.Range("C89").Value = "{0}, your certificate {1} has been successfully registered!"
$ws.Range("C89").Font.Size = 11

$ws.Range("A90").Value = "notification_certificate_issued_body"
$ws.Range("B90").Value = "{0}, podrás utilizar tú certificado proporcionando su identificador de emisión {1}"
$ws.Range("C90").Value = "{0}, you can use your certificate by providing its issue identifier {1}"
$ws.Range("C90").Font.Size = 11

$ws.Range("A91").Value = "notification_certificate_renewed_title"
$ws.Range("B91").Value = "{0}, Tú certificado {1} en Trust Certification System ha sido renovado!"
$ws.Range("C91").Value = "{0}, Your certificate {1} in the Trust Certification System has been renewed!"
$ws.Range("C91").Font.Size = 11

$ws.Range("A92").Value = "notification_certificate_renewed_body"
$ws.Range("B92").Value = "{0}, Tú certificado {1} ha sido renovado, podrás seguir utilizándolo con normalidad."
$ws.Range("C92").Value = "{0}, Your certificate {1} has been renewed, you can continue to use it normally."
$ws.Range("C92").Font.Size = 11

$ws.Range("A93").Value = "notification_certificate_request_accepted_title"
$ws.Range("B93").Value = "{0}, Tú solicitud de emisión del certificado {1} ha sido aceptada!"
$ws.Range("C93").Value = "{0}, Your request to issue the certificate {1} has been accepted!"

$ws.Range("A94").Value = "notification_certificate_request_accepted_body"
$ws.Range("B94").Value = "{0}, Tú solicitud de emisión del certificado {1} ha sido aceptada, en breve tú certificado será generado y almacenado en TCS."
$ws.Range("C94").Value = "{0}, Your request to issue the certificate {1} has been accepted, shortly your certificate will be generated and stored in TCS."

$ws.Range("A95").Value = "notification_certificate_request_rejected_title"
$ws.Range("B95").Value = "{0}, Tú solicitud de emisión de certificado {1} ha sido rechazada"
$ws.Range("C95").Value = "{0}, Your request to issue a certificate {1} has been rejected"

$ws.Range("A96").Value = "notification_certificate_request_rejected_body"
$ws.Range("B96").Value = "{0}, Tú solicitud de emisión del certificado {1} ha sido rechazada."
$ws.Range("C96").Value = "{0}, Your request to issue the certificate {1} has been rejected."

$ws.Range("A97").Value = "notification_certificate_visibility_changed_title"
$ws.Range("B97").Value = "{0}, La visibilidad de tú certificado {1} ha sido actualizada!."
$ws.Range("C97").Value = "{0}, The visibility of your certificate {1} has been updated!."

$ws.Range("A98").Value = "notification_certificate_changed_to_invisible_body"
$ws.Range("B98").Value = "{0}, La visibilidad de tú certificado {1} ha sido actualizada, tú certificado ahora no es visible."
$ws.Range("C98").Value = "{0}, The visibility of your certificate {1} has been updated, your certificate is now not visible."

$ws.Range("A99").Value = "notification_certificate_changed_to_visible_body"
$ws.Range("B99").Value = "{0}, La visibilidad de tú certificado {1} ha sido actualizada, tú certificado ahora es visible."
$ws.Range("C99").Value = "{0}, The visibility of your certificate {1} has been updated, your certificate is now visible."

$ws.Range("A100").Value = "notification_student_issue_certificate_requested_title"
$ws.Range("B100").Value = "{0}, tú solicitud de emisión de certificado ha sido registrada!"
$ws.Range("C100").Value = "{0}, your certificate issuance request has been registered!"

$ws.Range("A101").Value = "notification_student_issue_certificate_requested_body"
$ws.Range("B101").Value = "{0}, En breve recibirás información sobre el progreso de tú solicitud de emisión de certificado."
$ws.Range("C101").Value = "{0}, You will shortly receive information on the progress of your certificate issuance request"

$ws.Range("A102").Value = "notification_ca_issue_certificate_requested_title"
$ws.Range("B102").Value = "{0}, hay una nueva solicitud de emisión de certificado"
$ws.Range("C102").Value = "{0}, there is a new certificate issuance request."

$ws.Range("A103").Value = "notification_ca_issue_certificate_requested_body"
$ws.Range("B103").Value = "{0}, revisa la solicitud de emisión de certificado para completar el proceso."
$ws.Range("C103").Value = "{0}, review the certificate issuance request to complete the process."

# Update view selection to match final state
$ws.Range("C90").Select()